# "Changed nav bar to horizontal in a grid"
# - Grid_Mobile gets a new "Nav" row (row 3) describing a 3-item horizontal nav.
# - A new sheet "Grid_TabletPortrait" is added, a wider 8-column version of the
#   grid with its own (longer) horizontal Nav row.

$wb = $excel.ActiveWorkbook
$mobile = $wb.Worksheets.Item("Grid_Mobile")

# --- 1) Grid_Mobile: widen data columns B:J and add the Nav row -----------

$mobile.Range("B1:J1").ColumnWidth = 17.8

$mobile.Range("A3").Value = "Nav"
$mobile.Range("C3").Value = "logo (left aligned)"
$mobile.Range("E3").Value = "about (centered)"
$mobile.Range("G3").Value = "work (centered)"

# --- 2) Create Grid_TabletPortrait as a copy of Grid_Mobile, then rebuild --

$mobile.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$tablet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tablet.Name = "Grid_TabletPortrait"

# Clear out the copied Grid_Mobile content/selection so we can rebuild the
# wider 8-column layout from scratch.
$tablet.Cells.ClearContents()
$tablet.Range("B1:J1").ColumnWidth = 17.8
$tablet.Range("B1:R1").ColumnWidth = 14.5

# Row 1: Margin / Column N / Gutter headers across 8 columns
$tablet.Range("B1").Value = "Margin"
$tablet.Range("C1").Value = "Column 1"
$tablet.Range("D1").Value = "Gutter"
$tablet.Range("E1").Value = "Column 2"
$tablet.Range("F1").Value = "Gutter"
$tablet.Range("G1").Value = "Column 3"
$tablet.Range("H1").Value = "Gutter"
$tablet.Range("I1").Value = "Column 4"
$tablet.Range("J1").Value = "Gutter"
$tablet.Range("K1").Value = "Column 5"
$tablet.Range("L1").Value = "Gutter"
$tablet.Range("M1").Value = "Column 6"
$tablet.Range("N1").Value = "Gutter"
$tablet.Range("O1").Value = "Column 7"
$tablet.Range("P1").Value = "Gutter"
$tablet.Range("Q1").Value = "Column 8"
$tablet.Range("R1").Value = "Margin"

# Row 2: Dimensions
$tablet.Range("A2").Value = "Dimensions"
$tablet.Range("B2").Value = "24px"
$tablet.Range("C2").Value = "auto"
$tablet.Range("D2").Value = "16px"
$tablet.Range("E2").Value = "auto"
$tablet.Range("F2").Value = "16px"
$tablet.Range("G2").Value = "auto"
$tablet.Range("H2").Value = "16px"
$tablet.Range("I2").Value = "auto"
$tablet.Range("J2").Value = "16px"
$tablet.Range("K2").Value = "auto"
$tablet.Range("L2").Value = "16px"
$tablet.Range("M2").Value = "auto"
$tablet.Range("N2").Value = "16px"
$tablet.Range("O2").Value = "auto"
$tablet.Range("P2").Value = "16px"
$tablet.Range("Q2").Value = "auto"
$tablet.Range("R2").Value = "24px"

# Row 3: Nav
$tablet.Range("A3").Value = "Nav"
$tablet.Range("C3").Value = "Logo (left aligned)"
$tablet.Range("E3").Value = "About (centered)"
$tablet.Range("G3").Value = "Work (centered)"
$tablet.Range("I3").Value = "Education (centered)"
$tablet.Range("K3").Value = "Contact (centered)"
$tablet.Range("M3").Value = "."
$tablet.Range("O3").Value = "."
$tablet.Range("Q3").Value = "Jordan Tranchina (right-aligned)"
$tablet.Range("Q3").HorizontalAlignment = -4152

$tablet.Range("C3:G3").Select()
$mobile.Activate()
$mobile.Range("C3:G3").Select()
$tablet.Activate()
$tablet.Range("F7").Select()
